$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Range("H$row").Value = 2147482232
    $ws.Range("I$row").Value = 2147482236
    $ws.Range("AF$row").Value = 2147482228
}
